# Daily attendance processing
# For every row in the "Recorded By" column (G), rotate the comma-separated
# list of recorders one position to the left (the first name moves to the
# end of the list). Lists with a single entry are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -like "*,*") {
        $parts = $text -split ", "
        $count = $parts.Count

        $rotated = @()
        for ($i = 1; $i -lt $count; $i++) {
            $rotated += $parts[$i]
        }
        $rotated += $parts[0]

        $cell.Value = $rotated -join ", "
    }
}
